$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.084.74'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '3.247.96'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.97'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '185.03'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '3.808.86'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.89'
$ws.Range('E14').Value = '  -2.91%  '
$ws.Range('D15').Value = '68.074.50'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '3.243.17'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.47'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '396.83'
$ws.Range('E20').Value = '  +4.56%  '
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.40'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.186'
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.65'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.61'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.80'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.26'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  -4.78%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '161.81'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  +3.27%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '26.68'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.813'
$ws.Range('E39').Value = '  -3.12%  '
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.51'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('E42').Value = '  -3.94%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '41.22'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '25.36'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0685'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '2.612.21'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '336.12'
$ws.Range('E47').Value = '  -3.00%  '
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.32'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '31.15'
$ws.Range('E51').Value = '  +2.58%  '
